# Daily attendance processing - 2026-01-20 20:28:45
# Normalizes the "Recorded By" (column G) entries so that the primary
# recorder is listed first, swapping the order of the first two
# comma-separated recorders wherever needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact "Recorded By" values that need their first two entries swapped.
$targets = @(
    "System, backup@backdoor.com",
    "System, backup@backdoor.com, system",
    "dnasr281@gmail.com, System",
    "dnasr281@gmail.com, admin@admin.com"
)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$changed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq "") {
        continue
    }

    foreach ($target in $targets) {
        if ($val -eq $target) {
            $parts = $val -split ", "
            $first = $parts[0]
            $second = $parts[1]
            $rest = ""
            if ($parts.Count -gt 2) {
                for ($i = 2; $i -lt $parts.Count; $i++) {
                    $rest = $rest + ", " + $parts[$i]
                }
            }
            $newVal = $second + ", " + $first + $rest
            $cell.Value = $newVal
            $changed = $changed + 1
        }
    }
}

Write-Host "Recorded By values swapped:" $changed
